$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank columns at P:Q (shifts old P -> R)
$ws.Range("P1:Q1").EntireColumn.Insert()

# Headers
$ws.Range("P1").Value = "L1POST_drtn"
$ws.Range("Q1").Value = "L1PRE_mom_birth_country"

# Data for new columns P (L1POST_drtn) and Q (L1PRE_mom_birth_country)
$ws.Range("P2").Value = 11
$ws.Range("Q2").Value = 326.90804930000002
$ws.Range("P3").Value = 17
$ws.Range("Q3").Value = 401.80829310000001
$ws.Range("P4").Value = 14
$ws.Range("Q4").Value = 205.70598079999999
$ws.Range("P5").Value = 16
$ws.Range("Q5").Value = 287.37408110000001
$ws.Range("P6").Value = 10
$ws.Range("Q6").Value = 264.31338720000002
$ws.Range("P7").Value = 9
$ws.Range("Q7").Value = 160.0729868
$ws.Range("P8").Value = 4
$ws.Range("Q8").Value = 513.72649990000002
$ws.Range("P9").Value = 13
$ws.Range("Q9").Value = 228.18582169999999
$ws.Range("P10").Value = 4
$ws.Range("Q10").Value = 317.95951000000002
$ws.Range("P11").Value = 16
$ws.Range("Q11").Value = 139.98190360000001
$ws.Range("P12").Value = 21
$ws.Range("Q12").Value = 201.1972853
$ws.Range("P13").Value = 9
$ws.Range("Q13").Value = 256.03320100000002
$ws.Range("P14").Value = 9
$ws.Range("Q14").Value = 299.19605239999999
$ws.Range("P15").Value = 17
$ws.Range("Q15").Value = 388.51304010000001
$ws.Range("P16").Value = 12
$ws.Range("Q16").Value = 300.2996483
$ws.Range("P17").Value = 9
$ws.Range("Q17").Value = 323.79267199999998
$ws.Range("P18").Value = 13
$ws.Range("Q18").Value = 177.3339599
$ws.Range("P19").Value = 16
$ws.Range("Q19").Value = 311.07024259999997
$ws.Range("P20").Value = 13
$ws.Range("Q20").Value = 289.80873450000001
$ws.Range("P21").Value = 5
$ws.Range("Q21").Value = 402.90408289999999
$ws.Range("P22").Value = 6
$ws.Range("Q22").Value = 228.15228830000001
$ws.Range("P23").Value = 14
$ws.Range("Q23").Value = 288.13463159999998
$ws.Range("P24").Value = 9
$ws.Range("Q24").Value = 315.34628229999998
$ws.Range("P25").Value = 8
$ws.Range("Q25").Value = 316.90816610000002
$ws.Range("P26").Value = 7
$ws.Range("Q26").Value = 202.1801413
$ws.Range("P27").Value = 9
$ws.Range("Q27").Value = 238.10449389999999
$ws.Range("P28").Value = 7
$ws.Range("Q28").Value = 421.79011930000001
$ws.Range("P29").Value = 3
$ws.Range("Q29").Value = 58.705479590000003
$ws.Range("P30").Value = 19
$ws.Range("Q30").Value = 333.94450110000002
$ws.Range("P31").Value = 2
$ws.Range("Q31").Value = 212.13203440000001
$ws.Range("P32").Value = 12
$ws.Range("Q32").Value = 302.6422531
$ws.Range("P33").Value = 5
$ws.Range("Q33").Value = 432.00891189999999
$ws.Range("P34").Value = 11
$ws.Range("Q34").Value = 385.56199830000003
$ws.Range("P35").Value = 5
$ws.Range("Q35").Value = 282.60714780000001
$ws.Range("P36").Value = 4
$ws.Range("Q36").Value = 270.87312530000003
$ws.Range("P37").Value = 6
$ws.Range("Q37").Value = 335.51512630000002
$ws.Range("P38").Value = 7
$ws.Range("Q38").Value = 181.94949589999999
$ws.Range("P39").Value = 11
$ws.Range("Q39").Value = 181.41379119999999
$ws.Range("P40").Value = 6
$ws.Range("Q40").Value = 392.5803952
$ws.Range("P41").Value = 3
$ws.Range("Q41").Value = 234.29966569999999
$ws.Range("P42").Value = 10
$ws.Range("Q42").Value = 222.09409919999999
$ws.Range("P43").Value = 10
$ws.Range("Q43").Value = 252.4972167

# Data for new column S (L1PRE_mom_birth_country values), written after P/Q so
# new shared strings are appended in the same order as the source workbook
$ws.Range("S2").Value = "mexico"
$ws.Range("S3").Value = "el salvador"
$ws.Range("S4").Value = "el salvador"
$ws.Range("S5").Value = "mexico"
$ws.Range("S6").Value = "mexico"
$ws.Range("S7").Value = "usa"
$ws.Range("S8").Value = "mexico"
$ws.Range("S9").Value = "mexico"
$ws.Range("S10").Value = "mexico"
$ws.Range("S11").Value = "mexico"
$ws.Range("S12").Value = "mexico"
$ws.Range("S13").Value = "mexico"
$ws.Range("S14").Value = "mexico"
$ws.Range("S15").Value = "mexico"
$ws.Range("S16").Value = "mexico"
$ws.Range("S17").Value = "mexico"
$ws.Range("S18").Value = "usa"
$ws.Range("S19").Value = "honduras"
$ws.Range("S20").Value = "mexico"
$ws.Range("S21").Value = "el salvador"
$ws.Range("S22").Value = "mexico"
$ws.Range("S23").Value = "mexico"
$ws.Range("S24").Value = "usa"
$ws.Range("S25").Value = "mexico"
$ws.Range("S26").Value = "mexico"
$ws.Range("S27").Value = "guatemala"
$ws.Range("S28").Value = "mexico"
$ws.Range("S29").Value = "mexico"
$ws.Range("S30").Value = "mexico"
$ws.Range("S31").Value = "mexico"
$ws.Range("S32").Value = "mexico"
$ws.Range("S33").Value = "mexico"
$ws.Range("S34").Value = "mexico"
$ws.Range("S35").Value = "usa"
$ws.Range("S36").Value = "mexico"
$ws.Range("S37").Value = "mexico"
$ws.Range("S38").Value = "mexico"
$ws.Range("S39").Value = "mexico"
$ws.Range("S40").Value = "mexico"
$ws.Range("S41").Value = "mexico"
$ws.Range("S42").Value = "mexico"
$ws.Range("S43").Value = "mexico"

# Header for column S, written last (matches original authoring order)
$ws.Range("S1").Value = "L1POST_drtsd"

# Rows 41-43 have a custom row format (style 3); S-column cells in those rows
# should not inherit that formatting, matching the source workbook.
$ws.Range("S41").Style = "Normal"
$ws.Range("S42").Style = "Normal"
$ws.Range("S43").Style = "Normal"

# Column widths (character units); the stored OOXML width equals ColumnWidth + 5/6
$ws.Columns.Item(15).ColumnWidth = 15.666666666666666
$ws.Columns.Item(16).ColumnWidth = 13.998697916666666
$ws.Columns.Item(17).ColumnWidth = 13.998697916666666
$ws.Columns.Item(18).ColumnWidth = 15.498697916666666

# Selection
$ws.Range("N35").Select()
